$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second PC's clock speed replaces the placeholder header in D1
$ws.Range("D1").Value = "1.6 Ghz"

# Fill in the "ticks" results for the second PC (column E) for every
# treatment row that already has data in columns A/D/G.
$ws.Range("E3").Value = 18151
$ws.Range("E4").Value = 67
$ws.Range("E5").Value = 12
$ws.Range("E6").Value = 470
$ws.Range("E7").Value = 434
$ws.Range("E8").Value = 382
$ws.Range("E9").Value = 26654
$ws.Range("E10").Value = 51865
$ws.Range("E11").Value = 33492
$ws.Range("E12").Value = 8854
$ws.Range("E13").Value = 66
$ws.Range("E14").Value = 64
$ws.Range("E15").Value = 4147
$ws.Range("E16").Value = 5769
$ws.Range("E17").Value = 4252
$ws.Range("E18").Value = 672109
$ws.Range("E19").Value = 622476
$ws.Range("E20").Value = 363866

# Leave the selection where the author left it after data entry
$null = $ws.Range("H3").Select()
